$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it is used ---
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Latest Handback DateTime (column H) ---
# zh-cn keeps referencing the same shared string that "0001-01-01 00:00:00" used to use,
# whose text now becomes "2016-03-22 06:53:59".
$wsZhCn.Range("H2").Value = "2016-03-22 06:53:59"
$wsZhCn.Range("H3").Value = "2016-03-22 06:53:59"

# de-de moves to a distinct new datetime value.
$wsDeDe.Range("H2").Value = "2016-03-22 06:54:13"
$wsDeDe.Range("H3").Value = "2016-03-22 06:54:13"

# --- New columns F (Latest Target File) and G (Latest Handback File) ---
$mdName = "1761058d-58d9-4c45-a20e-e70b262a33d9.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fdcf8179c0b180f80c556d1b3494ef3c525e5c53/e2e/1761058d-58d9-4c45-a20e-e70b262a33d9.md"

$zhCnXlfName = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25bb827ec779a78c5325a01e8ae547469650ce62/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.zh-cn.xlf"

$deDeXlfName = "1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56abd67a8b8fdabb4e80ed514086018986e412e2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/1761058d-58d9-4c45-a20e-e70b262a33d9.07e8a009d35b56b1d145fb545e9bffebf8cf6370.de-de.xlf"

function Set-HandbackLinks($ws, $xlfName, $xlfUrl) {
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, $null, $null, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, $null, $null, $xlfName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, $null, $null, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, $null, $null, $xlfName) | Out-Null

    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276
    $ws.Range("G2").Font.Underline = 2
    $ws.Range("G2").Font.Color = 15570276
    $ws.Range("F3").Font.Underline = 2
    $ws.Range("F3").Font.Color = 15570276
    $ws.Range("G3").Font.Underline = 2
    $ws.Range("G3").Font.Color = 15570276
}

Set-HandbackLinks $wsZhCn $zhCnXlfName $zhCnXlfUrl
Set-HandbackLinks $wsDeDe $deDeXlfName $deDeXlfUrl

Write-Host "done"
